# feat: add 2022-Q1 data
#
# The previous "总计" (totals) sheet is renamed to "2022-Q1" and repurposed
# to hold that quarter's per-fund holdings table (same shape as the other
# quarterly sheets). A brand-new "总计" sheet is appended at the end,
# containing the same rolling summary table as before plus a new first row
# for 2022-Q1.

$wb = $excel.ActiveWorkbook

# A sheet that already has the "header row" / "index column" formatting we
# want to reuse (bold, centered, thin-bordered cells) so the new sheets
# pick up the same style instead of inventing new ones.
$styleSource = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1) Turn the existing "总计" sheet into the new "2022-Q1" detail sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q1.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$styleSource.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# code, name, scale(text), position(text), ratio(text), value(text, $null => numeric 0), rank(number)
$rows = @(
    @("012368", "摩根士丹利华鑫优享臻选六个月持有期混合型证券投资基金A", "5.76", "93.78", "5.55", "0.3197", 7),
    @("000309", "大摩品质生活精选股票", "4.36", "94.17", "5.45", "0.2376", 7),
    @("233006", "大摩领先优势混合", "4.12", "94.42", "5.29", "0.2179", 6),
    @("010322", "摩根士丹利华鑫新兴产业股票", "2.41", "94.11", "6.20", "0.1494", 6),
    @("012369", "摩根士丹利华鑫优享臻选六个月持有期混合型证券投资基金C", "0.40", "93.78", "5.55", "0.0222", 7),
    @("013072", "泰信医疗服务混合A", "0.10", "73.32", "3.37", "0.0034", 9),
    @("013073", "泰信医疗服务混合C", "0.00", "73.32", "3.37", $null, 9)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]

    $q1.Cells.Item($rowNum, 1).Value = $r

    $codeCell = $q1.Cells.Item($rowNum, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $data[0]

    $q1.Cells.Item($rowNum, 3).Value = $data[1]

    $scaleCell = $q1.Cells.Item($rowNum, 4)
    $scaleCell.NumberFormat = "@"
    $scaleCell.Value = $data[2]

    $posCell = $q1.Cells.Item($rowNum, 5)
    $posCell.NumberFormat = "@"
    $posCell.Value = $data[3]

    $ratioCell = $q1.Cells.Item($rowNum, 6)
    $ratioCell.NumberFormat = "@"
    $ratioCell.Value = $data[4]

    $valueCell = $q1.Cells.Item($rowNum, 7)
    if ($null -eq $data[5]) {
        $valueCell.Value = 0
    } else {
        $valueCell.NumberFormat = "@"
        $valueCell.Value = $data[5]
    }

    $q1.Cells.Item($rowNum, 8).Value = $data[6]
}

$styleSource.Range("A2:A8").Copy()
$q1.Range("A2:A8").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Append a fresh "总计" sheet after "2022-Q1" with the rolling summary.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($i = 0; $i -lt $totalHeaders.Length; $i++) {
    $total.Cells.Item(1, $i + 2).Value = $totalHeaders[$i]
}
$styleSource.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$totalRows = @(
    @("2022-Q1", 7, 0.95),
    @("2021-Q4", 11, 5.41),
    @("2021-Q3", 5, 0.91),
    @("2021-Q2", 9, 5.88),
    @("2021-Q1", 15, 6.33),
    @("2020-Q4", 7, 2.77)
)

for ($r = 0; $r -lt $totalRows.Length; $r++) {
    $rowNum = $r + 2
    $data = $totalRows[$r]

    $total.Cells.Item($rowNum, 1).Value = $r
    $total.Cells.Item($rowNum, 2).Value = $data[0]
    $total.Cells.Item($rowNum, 3).Value = $data[1]
    $total.Cells.Item($rowNum, 4).Value = $data[2]
}

$styleSource.Range("A2:A7").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
